$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Test"
$ws.Range("B1").Value = "Test"
$ws.Range("A2").Value = "Test"

$ws.Range("B2").Select()
